$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.846.61'
$ws.Range("E2").Value = '  -2.80%  '

$ws.Range("D3").Value = '3.353.40'
$ws.Range("E3").Value = '  -2.43%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'566.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.09%  '

$ws.Range("D6").Value = "'146.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("E8").Value = '  +0.32%  '

$ws.Range("D9").Value = "'7.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("E10").Value = '  -1.13%  '

$ws.Range("E11").Value = '  +2.16%  '

$ws.Range("D12").Value = '3.932.24'
$ws.Range("E12").Value = '  -2.32%  '

$ws.Range("D14").Value = "'27.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '

$ws.Range("D15").Value = '3.361.54'
$ws.Range("E15").Value = '  -2.51%  '

$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").Value = '60.884.00'
$ws.Range("E17").Value = '  -2.84%  '

$ws.Range("D18").Value = "'6.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("E19").Value = '  -0.79%  '

$ws.Range("D20").Value = "'8.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("D21").Value = "'376.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.89%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").Value = "'74.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").Value = '3.503.78'
$ws.Range("E25").Value = '  -2.27%  '

$ws.Range("E26").Value = '  -6.06%  '

$ws.Range("E27").Value = '  -3.84%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("E29").Value = '  -2.71%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("E32").Value = '  -3.51%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = "'1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = "'22.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.47%  '

$ws.Range("D35").Value = "'5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.58%  '

$ws.Range("D36").Value = "'168.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("E37").Value = '  -4.04%  '

$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("D39").Value = "'28.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.80%  '

$ws.Range("D40").Value = '3.391.45'
$ws.Range("E40").Value = '  -2.32%  '

$ws.Range("E41").Value = '  -2.79%  '

$ws.Range("E42").Value = '  -3.27%  '

$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = "'1.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.09%  '

$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = "'1.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.17%  '

$ws.Range("D46").Value = '2.494.47'
$ws.Range("E46").Value = '  -2.37%  '

$ws.Range("D47").Value = "'22.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.29%  '

$ws.Range("D48").Value = "'6.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.26%  '

$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("E50").Value = '  -2.04%  '

$ws.Range("E51").Value = '  +0.34%  '
